{"js": "// Apply the \"MY Address Assignment\" update:\n//  1. Merge the split \"Controller Triangle <color>\" runs for Blue/Green into\n//     single runs, give Yellow its own run (with the preceding <w:br/> in the\n//     same run), and append three new controller lines (Purple 0x3D,\n//     Black 0x3E, Orange 0x3F).\n//  2. Move the \"_GoBack\" bookmark from the end of the first paragraph to the\n//     end of the (now longer) second paragraph.\n//  3. Move the <w:lastRenderedPageBreak/> marker from the \"Router 1\" table\n//     cell to the \"Router Name\" header cell (the content reflow pushed the\n//     rendered page break earlier).\n\nfunction wrapBody(inner) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body>\" +\n    inner +\n    \"</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\n// --- Remove the stale \"_GoBack\" bookmark at the end of the first paragraph ---\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- Rebuild the second paragraph (Manager / Controller Triangle list) ---\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst controllerParagraph = paragraphs.items[1];\nconst controllerRange = controllerParagraph.getRange(Word.RangeLocation.whole);\n\nconst newParagraphXml = wrapBody(\n  \"<w:p>\" +\n    \"<w:r><w:t>Manager MY: 0x01</w:t></w:r>\" +\n    \"<w:r><w:br/></w:r>\" +\n    \"<w:r><w:br/></w:r>\" +\n    \"<w:r><w:t>Controller Triangle Red MY: 0x39</w:t></w:r>\" +\n    \"<w:r><w:br/><w:t>Controller Triangle Blue MY: 0x3A</w:t></w:r>\" +\n    \"<w:r><w:br/><w:t>Controller Triangle Green MY: 0x3B</w:t></w:r>\" +\n    \"<w:r><w:br/><w:t>Controller Triangle Yellow MY: 0x3C</w:t></w:r>\" +\n    \"<w:r><w:br/><w:t>Controller Triangle Purple MY: 0x3D</w:t></w:r>\" +\n    \"<w:r><w:br/><w:t>Controller Triangle Black MY: 0x3E</w:t></w:r>\" +\n    \"<w:r><w:br/><w:t>Controller Triangle Orange MY: 0x3F</w:t></w:r>\" +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n    \"</w:p>\"\n);\n\ncontrollerRange.insertOoxml(newParagraphXml, Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Move <w:lastRenderedPageBreak/> from \"Router 1\" to \"Router Name\" ---\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst routerTable = tables.items[1];\n\nconst headerCell = routerTable.getCell(0, 0);\nconst headerParagraph = headerCell.body.paragraphs.getFirst();\nconst headerRange = headerParagraph.getRange(Word.RangeLocation.whole);\nconst headerXml = wrapBody(\n  \"<w:p>\" +\n    '<w:pPr><w:jc w:val=\"center\"/><w:rPr><w:u w:val=\"single\"/></w:rPr></w:pPr>' +\n    '<w:r><w:rPr><w:u w:val=\"single\"/></w:rPr><w:lastRenderedPageBreak/><w:t>Router Name</w:t></w:r>' +\n    \"</w:p>\"\n);\nheaderRange.insertOoxml(headerXml, Word.InsertLocation.replace);\nawait context.sync();\n\nconst router1Cell = routerTable.getCell(2, 0);\nconst router1Paragraph = router1Cell.body.paragraphs.getFirst();\nconst router1Range = router1Paragraph.getRange(Word.RangeLocation.whole);\nconst router1Xml = wrapBody(\n  \"<w:p>\" +\n    '<w:pPr><w:jc w:val=\"center\"/></w:pPr>' +\n    \"<w:r><w:t>Router 1</w:t></w:r>\" +\n    \"</w:p>\"\n);\nrouter1Range.insertOoxml(router1Xml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Apply the \"MY Address Assignment\" update:\n#  1. Merge the split \"Controller Triangle <color>\" runs for Blue/Green into\n#     single runs, give Yellow its own run (with the preceding <w:br/> in the\n#     same run), and append three new controller lines (Purple 0x3D,\n#     Black 0x3E, Orange 0x3F).\n#  2. Move the \"_GoBack\" bookmark from the end of the first paragraph to the\n#     end of the (now longer) second paragraph.\n#  3. Move the <w:lastRenderedPageBreak/> marker from the \"Router 1\" table\n#     cell to the \"Router Name\" header cell (the content reflow pushed the\n#     rendered page break earlier).\n\n$d = $word.ActiveDocument\n\n# --- Remove the stale \"_GoBack\" bookmark at the end of the first paragraph ---\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# --- Rebuild the second paragraph (Manager / Controller Triangle list) ---\n$controllerParagraph = $d.Paragraphs.Item(2)\n$controllerRange = $controllerParagraph.Range\n$newParagraphXml = @\"\n<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>\n<w:r><w:t>Manager MY: 0x01</w:t></w:r>\n<w:r><w:br/></w:r>\n<w:r><w:br/></w:r>\n<w:r><w:t>Controller Triangle Red MY: 0x39</w:t></w:r>\n<w:r><w:br/><w:t>Controller Triangle Blue MY: 0x3A</w:t></w:r>\n<w:r><w:br/><w:t>Controller Triangle Green MY: 0x3B</w:t></w:r>\n<w:r><w:br/><w:t>Controller Triangle Yellow MY: 0x3C</w:t></w:r>\n<w:r><w:br/><w:t>Controller Triangle Purple MY: 0x3D</w:t></w:r>\n<w:r><w:br/><w:t>Controller Triangle Black MY: 0x3E</w:t></w:r>\n<w:r><w:br/><w:t>Controller Triangle Orange MY: 0x3F</w:t></w:r>\n<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n<w:bookmarkEnd w:id=\"0\"/>\n</w:p>\n\"@\n$controllerRange.InsertXML($newParagraphXml)\n\n# --- Move <w:lastRenderedPageBreak/> from \"Router 1\" to \"Router Name\" ---\n$routerTable = $d.Tables.Item(2)\n\n$headerCell = $routerTable.Cell(1, 1)\n$headerParagraph = $headerCell.Range.Paragraphs.Item(1)\n$headerRange = $headerParagraph.Range\n$headerXml = @\"\n<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>\n<w:pPr><w:jc w:val=\"center\"/><w:rPr><w:u w:val=\"single\"/></w:rPr></w:pPr>\n<w:r><w:rPr><w:u w:val=\"single\"/></w:rPr><w:lastRenderedPageBreak/><w:t>Router Name</w:t></w:r>\n</w:p>\n\"@\n$headerRange.InsertXML($headerXml)\n\n$router1Cell = $routerTable.Cell(3, 1)\n$router1Paragraph = $router1Cell.Range.Paragraphs.Item(1)\n$router1Range = $router1Paragraph.Range\n$router1Xml = @\"\n<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>\n<w:pPr><w:jc w:val=\"center\"/></w:pPr>\n<w:r><w:t>Router 1</w:t></w:r>\n</w:p>\n\"@\n$router1Range.InsertXML($router1Xml)\n"}
